$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1035.16
$ws.Range("J17").Value = 1057.4584
$ws.Range("L17").Value = 3172.3752
$ws.Range("N17").Value = -3508.3752

$ws.Range("H86").Value = 42749.9
$ws.Range("I86").Value = 64331.156
$ws.Range("J86").Value = 1745.5
$ws.Range("K86").Value = 64331.156
$ws.Range("L86").Value = 1745.5
$ws.Range("M86").Value = -63208.156
$ws.Range("N86").Value = -3991.5

$ws.Range("H89").Value = 42749.9
$ws.Range("I89").Value = 64331.156
$ws.Range("J89").Value = 1745.5
$ws.Range("K89").Value = 321655.78
$ws.Range("L89").Value = 8727.5
$ws.Range("M89").Value = -316039.78
$ws.Range("N89").Value = -19959.5

$ws.Range("H127").Value = 1304.3
$ws.Range("I127").Value = 1182.1666
$ws.Range("K127").Value = 3546.4998
$ws.Range("M127").Value = 1413.5002

$ws.Range("H132").Value = 4713.2646
$ws.Range("I132").Value = 3681.8076
$ws.Range("J132").Value = 8065.5
$ws.Range("K132").Value = 11045.4228
$ws.Range("L132").Value = 24196.5
$ws.Range("M132").Value = -8515.4228
$ws.Range("N132").Value = -29256.5

$ws.Range("H138").Value = 1528.8933
$ws.Range("I138").Value = 1165.6792
$ws.Range("J138").Value = 2403.9092
$ws.Range("K138").Value = 3497.0376
$ws.Range("L138").Value = 7211.7276
$ws.Range("M138").Value = 1642.9624
$ws.Range("N138").Value = -17491.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3361.38
$ws.Range("I32").Value = 3365.0303
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 3365.0303
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -3078.0303
$ws.Range("N32").Value = -3574

$ws.Range("H88").Value = 1785.75
$ws.Range("I88").Value = 1725
$ws.Range("J88").Value = 1806
$ws.Range("K88").Value = 1725
$ws.Range("L88").Value = 1806
$ws.Range("M88").Value = -1319
$ws.Range("N88").Value = -2618

$ws.Range("H91").Value = 1785.75
$ws.Range("I91").Value = 1725
$ws.Range("J91").Value = 1806
$ws.Range("K91").Value = 1725
$ws.Range("L91").Value = 1806
$ws.Range("M91").Value = -321
$ws.Range("N91").Value = -4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1635.3334
$ws.Range("I86").Value = 1362.4
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1362.4
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -239.4000000000001
$ws.Range("N86").Value = -5246

$ws.Range("H89").Value = 1635.3334
$ws.Range("I89").Value = 1362.4
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 6812
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -1196
$ws.Range("N89").Value = -26232

$ws.Range("H134").Value = 1682.4459
$ws.Range("I134").Value = 921.0192
$ws.Range("J134").Value = 3482.182
$ws.Range("K134").Value = 2763.0576
$ws.Range("L134").Value = 10446.546
$ws.Range("M134").Value = -228.0576000000001
$ws.Range("N134").Value = -15516.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1880.28
$ws.Range("I31").Value = 1119.9482
$ws.Range("J31").Value = 2930.262
$ws.Range("K31").Value = 1119.9482
$ws.Range("L31").Value = 2930.262
$ws.Range("M31").Value = -824.9482
$ws.Range("N31").Value = -3520.262

$ws.Range("H34").Value = 1880.28
$ws.Range("I34").Value = 1119.9482
$ws.Range("J34").Value = 2930.262
$ws.Range("K34").Value = 1119.9482
$ws.Range("L34").Value = 2930.262
$ws.Range("M34").Value = -917.9482
$ws.Range("N34").Value = -3334.262

$ws.Range("H132").Value = 2237.138
$ws.Range("I132").Value = 1441.5143
$ws.Range("J132").Value = 3447.8696
$ws.Range("K132").Value = 4324.5429
$ws.Range("L132").Value = 10343.6088
$ws.Range("M132").Value = -1794.5429
$ws.Range("N132").Value = -15403.6088

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1802.375
$ws.Range("I5").Value = 1458.5
$ws.Range("J5").Value = 2146.25
$ws.Range("K5").Value = 4375.5
$ws.Range("L5").Value = 6438.75
$ws.Range("M5").Value = -4263.5
$ws.Range("N5").Value = -6662.75

$ws.Range("H122").Value = 1984.2549
$ws.Range("I122").Value = 449.8421
$ws.Range("J122").Value = 2895.3125
$ws.Range("K122").Value = 4048.5789
$ws.Range("L122").Value = 26057.8125
$ws.Range("M122").Value = -1598.5789
$ws.Range("N122").Value = -30957.8125

$ws.Range("H135").Value = 1802.375
$ws.Range("I135").Value = 1458.5
$ws.Range("J135").Value = 2146.25
$ws.Range("K135").Value = 13126.5
$ws.Range("L135").Value = 19316.25
$ws.Range("M135").Value = -10591.5
$ws.Range("N135").Value = -24386.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3297.1667
$ws.Range("I62").Value = 2718.4443
$ws.Range("J62").Value = 5033.3335
$ws.Range("K62").Value = 2718.4443
$ws.Range("L62").Value = 5033.3335
$ws.Range("M62").Value = -2094.4443
$ws.Range("N62").Value = -6281.3335

$ws.Range("H63").Value = 21120.834
$ws.Range("I63").Value = 2226
$ws.Range("J63").Value = 24899.8
$ws.Range("K63").Value = 2226
$ws.Range("L63").Value = 24899.8
$ws.Range("M63").Value = -1602
$ws.Range("N63").Value = -26147.8

$ws.Range("H65").Value = 3297.1667
$ws.Range("I65").Value = 2718.4443
$ws.Range("J65").Value = 5033.3335
$ws.Range("K65").Value = 13592.2215
$ws.Range("L65").Value = 25166.6675
$ws.Range("M65").Value = -10472.2215
$ws.Range("N65").Value = -31406.6675

$ws.Range("H66").Value = 21120.834
$ws.Range("I66").Value = 2226
$ws.Range("J66").Value = 24899.8
$ws.Range("K66").Value = 6678
$ws.Range("L66").Value = 74699.39999999999
$ws.Range("M66").Value = -3558
$ws.Range("N66").Value = -80939.39999999999

$ws.Range("H70").Value = 15000
$ws.Range("J70").Value = 15000
$ws.Range("L70").Value = 15000
$ws.Range("N70").Value = -15630

$ws.Range("H73").Value = 15000
$ws.Range("J73").Value = 15000
$ws.Range("L73").Value = 15000
$ws.Range("N73").Value = -17184

$ws.Range("H132").Value = 1662.8679
$ws.Range("I132").Value = 1424.6052
$ws.Range("J132").Value = 2266.4666
$ws.Range("K132").Value = 4273.8156
$ws.Range("L132").Value = 6799.399800000001
$ws.Range("M132").Value = -1743.8156
$ws.Range("N132").Value = -11859.3998
